# Working on session object meta data akquisition
# Remove the two "in progress" task rows (B9: Events/InitialDiagnosis note,
# B10: InitialDiagnosis nested details note) from the ToDo list, while
# keeping the row numbers of the remaining tasks unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 9 and 10 (this shifts everything below up), then immediately
# insert two blank rows back at the same position so the later rows (12,
# 14, 16, 20, 22, ...) keep their original row numbers while rows 9/10
# themselves no longer carry any content.
$ws.Rows("9:10").Delete()
$ws.Rows("9:10").Insert()

# Update the view: scroll position back to the top-left and move the
# active selection to C12.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C12").Select()
